$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.645.06'
$ws.Range("E2").Value = '  -4.37%  '
$ws.Range("D3").Value = '1.846.14'
$ws.Range("E3").Value = '  -3.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4247'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -7.00%  '
$ws.Range("E8").Value = '  -4.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.82'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07208'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8986'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.84%  '
$ws.Range("E12").Value = '  -7.15%  '
$ws.Range("D13").Value = '1.836.41'
$ws.Range("E13").Value = '  -5.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.574'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.344'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06797'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '77.52'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -8.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008852'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9999'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.69%  '
$ws.Range("D22").Value = '27.614.64'
$ws.Range("E22").Value = '  -4.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.956'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.69'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.24%  '
$ws.Range("D25").Value = '2.059.04'
$ws.Range("E25").Value = '  -4.19%  '
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.13%  '
$ws.Range("E28").Value = '  -4.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.324'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '111.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.768'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08901'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7752'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -10.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.501'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -11.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.847'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.081'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -12.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.000'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05440'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.094'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.979'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01925'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5044'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.775'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1633'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.74%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.259'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -11.29%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06619'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '106.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.66%  '
$ws.Range("E48").Value = '  -8.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9997'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.645'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.48%  '
